# Update rows 2-11 with the new TPM-derived values and re-labelled clusters
# (cluster "MuSCs" -> "Inflammatory-Mac" inserted as sender category; "MuSCs" dropped as a target)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Uts2"
$ws.Range("C2").Value = "Uts2r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1789903333333333
$ws.Range("H2").Value = 0.536971
$ws.Range("I2").Value = 0.1405940142989478
$ws.Range("J2").Value = 0.1620097168840439
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2516033333333333
$ws.Range("N2").Value = 0.75481
$ws.Range("O2").Value = 0.5477584147437079
$ws.Range("P2").Value = 0.5477584147437079
$ws.Range("Q2").Value = 0.04503456450111112
$ws.Range("R2").Value = 0.40531108051
$ws.Range("S2").Value = 0.07701155439484586
$ws.Range("T2").Value = 0.08874218569348084

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Uts2"
$ws.Range("C3").Value = "Uts2r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1789903333333333
$ws.Range("H3").Value = 0.536971
$ws.Range("I3").Value = 0.1405940142989478
$ws.Range("J3").Value = 0.1620097168840439
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2077293333333333
$ws.Range("N3").Value = 0.623188
$ws.Range("O3").Value = 0.4522415852562921
$ws.Range("P3").Value = 0.4522415852562921
$ws.Range("Q3").Value = 0.03718154261644444
$ws.Range("R3").Value = 0.3346338835479999
$ws.Range("S3").Value = 0.06358245990410197
$ws.Range("T3").Value = 0.0732675311905631

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Uts2"
$ws.Range("C4").Value = "Uts2r"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.297118
$ws.Range("H4").Value = 0.8913540000000001
$ws.Range("I4").Value = 0.2333813874891276
$ws.Range("J4").Value = 0.2689307414803781
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2516033333333333
$ws.Range("N4").Value = 0.75481
$ws.Range("O4").Value = 0.5477584147437079
$ws.Range("P4").Value = 0.5477584147437079
$ws.Range("Q4").Value = 0.07475587919333335
$ws.Range("R4").Value = 0.67280291274
$ws.Range("S4").Value = 0.1278366188417316
$ws.Range("T4").Value = 0.1473090766291418

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Uts2"
$ws.Range("C5").Value = "Uts2r"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.297118
$ws.Range("H5").Value = 0.8913540000000001
$ws.Range("I5").Value = 0.2333813874891276
$ws.Range("J5").Value = 0.2689307414803781
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2077293333333333
$ws.Range("N5").Value = 0.623188
$ws.Range("O5").Value = 0.4522415852562921
$ws.Range("P5").Value = 0.4522415852562921
$ws.Range("Q5").Value = 0.06172012406133334
$ws.Range("R5").Value = 0.555481116552
$ws.Range("S5").Value = 0.1055447686473961
$ws.Range("T5").Value = 0.1216216648512363

# Row 6: Inflammatory-Mac -> ECs
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Uts2"
$ws.Range("C6").Value = "Uts2r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03373033333333333
$ws.Range("H6").Value = 0.101191
$ws.Range("I6").Value = 0.02649463174161143
$ws.Range("J6").Value = 0.03053037363510002
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2516033333333333
$ws.Range("N6").Value = 0.75481
$ws.Range("O6").Value = 0.5477584147437079
$ws.Range("P6").Value = 0.5477584147437079
$ws.Range("Q6").Value = 0.008486664301111111
$ws.Range("R6").Value = 0.07637997871
$ws.Range("S6").Value = 0.0145126574820034
$ws.Range("T6").Value = 0.01672326906389548

# Row 7: Inflammatory-Mac -> FAPs
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Uts2"
$ws.Range("C7").Value = "Uts2r"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03373033333333333
$ws.Range("H7").Value = 0.101191
$ws.Range("I7").Value = 0.02649463174161143
$ws.Range("J7").Value = 0.03053037363510002
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2077293333333333
$ws.Range("N7").Value = 0.623188
$ws.Range("O7").Value = 0.4522415852562921
$ws.Range("P7").Value = 0.4522415852562921
$ws.Range("Q7").Value = 0.007006779656444444
$ws.Range("R7").Value = 0.06306101690799999
$ws.Range("S7").Value = 0.01198197425960803
$ws.Range("T7").Value = 0.01380710457120454

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Uts2"
$ws.Range("C8").Value = "Uts2r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.504865
$ws.Range("H8").Value = 1.00973
$ws.Range("I8").Value = 0.3965632987388795
$ws.Range("J8").Value = 0.304646007753353
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2516033333333333
$ws.Range("N8").Value = 0.75481
$ws.Range("O8").Value = 0.5477584147437079
$ws.Range("P8").Value = 0.5477584147437079
$ws.Range("Q8").Value = 0.1270257168833333
$ws.Range("R8").Value = 0.7621543013
$ws.Range("S8").Value = 0.2172208838627441
$ws.Range("T8").Value = 0.166872414264976

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Uts2"
$ws.Range("C9").Value = "Uts2r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.504865
$ws.Range("H9").Value = 1.00973
$ws.Range("I9").Value = 0.3965632987388795
$ws.Range("J9").Value = 0.304646007753353
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2077293333333333
$ws.Range("N9").Value = 0.623188
$ws.Range("O9").Value = 0.4522415852562921
$ws.Range("P9").Value = 0.4522415852562921
$ws.Range("Q9").Value = 0.1048752698733333
$ws.Range("R9").Value = 0.62925161924
$ws.Range("S9").Value = 0.1793424148761354
$ws.Range("T9").Value = 0.137773593488377

# Row 10: Resolving-Mac -> ECs
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Uts2"
$ws.Range("C10").Value = "Uts2r"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.258397
$ws.Range("H10").Value = 0.775191
$ws.Range("I10").Value = 0.2029666677314337
$ws.Range("J10").Value = 0.2338831602471249
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2516033333333333
$ws.Range("N10").Value = 0.75481
$ws.Range("O10").Value = 0.5477584147437079
$ws.Range("P10").Value = 0.5477584147437079
$ws.Range("Q10").Value = 0.06501354652333333
$ws.Range("R10").Value = 0.58512191871
$ws.Range("S10").Value = 0.111176700162383
$ws.Range("T10").Value = 0.1281114690922137

# Row 11: Resolving-Mac -> FAPs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Uts2"
$ws.Range("C11").Value = "Uts2r"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.258397
$ws.Range("H11").Value = 0.775191
$ws.Range("I11").Value = 0.2029666677314337
$ws.Range("J11").Value = 0.2338831602471249
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2077293333333333
$ws.Range("N11").Value = 0.623188
$ws.Range("O11").Value = 0.4522415852562921
$ws.Range("P11").Value = 0.4522415852562921
$ws.Range("Q11").Value = 0.05367663654533333
$ws.Range("R11").Value = 0.4830897289079999
$ws.Range("S11").Value = 0.09178996756905067
$ws.Range("T11").Value = 0.1057716911549112

# The new table only has 10 data rows (2-11); remove the former rows 12 and 13
$ws.Range("A12:A13").EntireRow.Delete()

Write-Host "Edit complete. UsedRange: $($ws.UsedRange.Address())"
